# Add three "Extra" columns (Extra1, Extra2, Extra3) to the end of the
# header row on each of the three worksheets, then leave the selection on
# each sheet positioned just past the new columns (mirrors what Excel
# leaves behind after typing headers across a row and tabbing through).
#
# Sheet layout (by name, since tab order == workbook.xml <sheets> order):
#   PROSPECT   (sheet1): headers A1:I1 -> add J1:L1   (Extra1, Extra2, Extra3)
#   STUDENTS   (sheet2): headers A1:P1 -> add Q1:S1   (Extra1, Extra2, Extra3)
#   EX-STUDENT (sheet3): headers A1:L1 -> add M1:O1   (Extra1, Extra2, Extra3)

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("PROSPECT")
$ws2 = $wb.Worksheets.Item("STUDENTS")
$ws3 = $wb.Worksheets.Item("EX-STUDENT")

# --- STUDENTS (sheet2) ---
$ws2.Range("Q1").Value = "Extra1"
$ws2.Range("R1").Value = "Extra2"
$ws2.Range("S1").Value = "Extra3"
$ws2.Range("Q1:S1048576").Select() | Out-Null

# --- EX-STUDENT (sheet3) ---
$ws3.Range("M1").Value = "Extra1"
$ws3.Range("N1").Value = "Extra2"
$ws3.Range("O1").Value = "Extra3"
$ws3.Range("M1:O1048576").Select() | Out-Null

# --- PROSPECT (sheet1) --- done last so the workbook's active tab stays
# on PROSPECT, matching the original file (tabSelected stays on sheet1).
$ws1.Range("J1").Value = "Extra1"
$ws1.Range("K1").Value = "Extra2"
$ws1.Range("L1").Value = "Extra3"
$ws1.Range("N21").Select() | Out-Null
